# "debug irregular date string"
# - PayrollSchedule row 4 (March 2025) becomes December 2025:
#     PayrollMonth 202503 -> 202512
#     CutoffDate   45741 (2025-03-25) -> 45996 (2025-12-05)
#     PayDate      45747 (2025-03-31) -> 46020 (2025-12-29)
# - Column C (PayDate) widens slightly to fit the new best-fit date text.
# - The active sheet/tab moves from ExchangeRates back to PayrollSchedule,
#   with cell B4 selected there.
# - The saved workbook window position moves too (xWindow/yWindow).

$wb = $excel.ActiveWorkbook

$payroll = $wb.Worksheets.Item("PayrollSchedule")

# Update the irregular row of data (PayrollMonth / CutoffDate / PayDate).
$payroll.Range("A4").Value = 202512
$payroll.Range("B4").Value = 45996
$payroll.Range("C4").Value = 46020

# Widen column C (PayDate) to its new best-fit-like width.
$payroll.Columns.Item(3).ColumnWidth = 10.37

# Move the window position (best effort - engine may not expose this).
try {
    $win = $excel.ActiveWindow
    $win.Left = 780
    $win.Top = 780
} catch {}

# Make PayrollSchedule the active sheet/tab again, with B4 selected,
# which also clears the previous tabSelected flag on ExchangeRates.
$payroll.Activate()
$payroll.Range("B4").Select()
